$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.384.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.307.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.532"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.72%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "51.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.668.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.306.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.810"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.283.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0928"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "242.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("E25").Value = "  +2.78%  "
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.68%  "
$ws.Range("E29").Value = "  +7.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("E36").Value = "  +6.19%  "
$ws.Range("E37").Value = "  +0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.07%  "
$ws.Range("E40").Value = "  +2.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.94%  "
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0297"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.986.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("E47").Value = "  +3.63%  "
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("E51").Value = "  +9.18%  "
